# "Update Data Sources from LFX"
#
# The deck's data tables (slides 9, 13, 14, 15, 16, 18) were all stamped with
# an old table-style GUID; refresh them to the new style GUID. Walk every
# slide/shape instead of hard-coding indices, and only touch tables that are
# still on the old style, so the script is robust to reordering and safe to
# re-run.

$p = $ppt.ActivePresentation

$oldStyleId = "{76AEB783-E90B-4FD3-9B3C-ECA70C141418}"
$newStyleId = "{8D3E67DC-7664-4D72-9394-26F2834BD407}"

for ($slideIndex = 1; $slideIndex -le $p.Slides.Count; $slideIndex++) {
    $slide = $p.Slides.Item($slideIndex)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)

        if ($shape.HasTable) {
            $table = $shape.Table

            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
